$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text in the source data
# (e.g. "600.75", "0.200", "3.416.73"). Excel auto-converts numeric-looking
# strings assigned via .Value into real numbers, which would corrupt values
# like trailing zeros or thousand-dot formatting. Force a Text format on each
# of these cells first so the subsequent assignment keeps them as text.
$priceCells = @(
    "D2", "D3", "D5", "D6", "D7", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D50"
)
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "70.916.52"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "3.598.51"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "600.75"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").Value = "174.13"
$ws.Range("E6").Value = "  +1.67%  "
$ws.Range("D7").Value = "3.593.19"
$ws.Range("E7").Value = "  +2.31%  "
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "0.200"
$ws.Range("E10").Value = "  +6.10%  "
$ws.Range("D11").Value = "7.46"
$ws.Range("E11").Value = "  +8.15%  "
$ws.Range("D12").Value = "0.592"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "47.13"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "0.0000280"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").Value = "4.174.48"
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("D16").Value = "8.47"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "616.28"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "3.601.50"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("D19").Value = "70.993.39"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").Value = "17.52"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").Value = "0.891"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").Value = "9.31"
$ws.Range("E23").Value = "  -16.42%  "
$ws.Range("D24").Value = "16.07"
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("D25").Value = "97.70"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "3.79"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "2.67"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "34.13"
$ws.Range("E29").Value = "  +4.24%  "
$ws.Range("D30").Value = "9.26"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "8.52"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").Value = "3.09"
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("D33").Value = "7.29"
$ws.Range("E33").Value = "  +4.87%  "
$ws.Range("D34").Value = "1.31"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").Value = "641.68"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "3.73"
$ws.Range("E36").Value = "  +6.82%  "
$ws.Range("D37").Value = "0.102"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").Value = "10.89"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").Value = "0.0480"
$ws.Range("E39").Value = "  +5.74%  "
$ws.Range("D40").Value = "57.43"
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "0.143"
$ws.Range("E42").Value = "  +5.73%  "
$ws.Range("D43").Value = "3.416.73"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").Value = "0.326"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").Value = "0.0₃0721"
$ws.Range("E45").Value = "  +3.02%  "
$ws.Range("D46").Value = "2.99"
$ws.Range("E46").Value = "  +8.50%  "
$ws.Range("D47").Value = "33.10"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").Value = "2.69"
$ws.Range("E48").Value = "  +5.92%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").Value = "132.94"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("E51").Value = "  -0.11%  "
